# Final-Project.docx — "big mid-semester update" edits.
#
# Applies the text-content changes from the commit:
#  1. Due-date paragraph: deadline/requirements text rewritten (new date,
#     Blackboard submission instructions), which also removes the stray
#     "_GoBack" bookmark that sat inside the old wording.
#  2. "cite your data source" bullet: add "(and share!)".
#  3. Kaggle bullet: "a lot more" -> "many more".
#  4. Text-analysis bullet: drop "maybe " before "constraining".
#  5. Tic-Tac-Toe bullet: "it plays at random" -> "the computer plays at
#     random" (disambiguating the pronoun).
#
# wdReplace constants used with Find.Execute's Replace argument:
#   wdReplaceNone = 0, wdReplaceOne = 1, wdReplaceAll = 2

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Output ("NOT FOUND: " + $find)
    }
}

# 1) Deadline paragraph — new requirement (uploaded to Blackboard twice) and
#    new due date/time; this span also swallows the old "_GoBack" bookmark.
Replace-Text `
    "Your project should be complete and production-ready by 6pm on May 6th (our assigned exam period): there should be good interaction" `
    "Your project should be complete, production-ready, and uploaded to Blackboard (twice—once as an assignment submission and once in the Discussion Board) by midnight on Sunday, December 13: there should be good interaction"

# 2) Data-source bullet — "keep (and share!) notes about..."
Replace-Text `
    " sure to cite your data source, and keep notes about any rows or columns you" `
    " sure to cite your data source, and keep (and share!) notes about any rows or columns you"

# 3) Kaggle bullet — "a lot more" -> "many more"
Replace-Text `
    " Kaggle’s got a lot more general data sets, and you’re welcome to go out and find or make your own data to analyze." `
    " Kaggle’s got many more general data sets, and you’re welcome to go out and find or make your own data to analyze."

# 4) Text-analysis bullet — drop "maybe "
Replace-Text `
    " recommend maybe constraining yourself to a single short story or long poem, rather than a whole book.)" `
    " recommend constraining yourself to a single short story or long poem, rather than a whole book.)"

# 5) Games bullet — clarify "the computer plays at random"
Replace-Text `
    " prepared to write. You could make a Sudoku solver, a smart Tic-Tac-Toe (first pass: it plays at random; improvement: it tries to win, after letting the player play first), 2048, or Battleship. Most of these don’t " `
    " prepared to write. You could make a Sudoku solver, a smart Tic-Tac-Toe (first pass: the computer plays at random; improvement: it tries to win, after letting the player play first), 2048, or Battleship. Most of these don’t "

Write-Output "done"
